$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("PDiCECpDoC")

# Update the header label on the "PDiCECpDoC" sheet to clarify it is dimensionless
$ws.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"

# Reflect the active selection being on B2 (the value cell) in that sheet
$ws.Range("B2").Select()

# Keep the "About" sheet as the active/selected tab, as in the source workbook
$wsAbout.Activate()
